$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "Planilha1" -> "5W2H"
$ws.Name = "5W2H"

# Update the active selection on the sheet to A12
$ws.Range("A12").Select()
